$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 ("质控组" roster) ---
# Remove the "张悦" row (all-zero row, originally row 2); remaining rows
# shift up so 冷雪 becomes row 2 and 屈昂 becomes row 3.
$ws1.Range("A2").EntireRow.Delete()

# The team name changed from "质控组" to "北京组" for the remaining rows.
$ws1.Range("A2").Value = "北京组"
$ws1.Range("A3").Value = "北京组"

# Selection on Sheet1 now spans the entire 2nd row.
[void]$ws1.Rows(2).Select()

# --- Sheet2 ("总体" summary) ---
# Same team-name rename on the summary sheet.
$ws2.Range("A2").Value = "北京组"

# Updated lims-error figure.
$ws2.Range("K2").Value = 4.66

# Sheet2's view had scrolled right (topLeftCell F1); reset the scroll
# position back to the top-left and move the selection to A2.
[void]$ws2.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
[void]$ws2.Range("A2").Select()

# Restore Sheet1 as the active/tab-selected sheet.
[void]$ws1.Activate()
